# Replace "OIE" with "WOAH" throughout the relevant text cells of the
# workbook (RVFV story map data contents), matching the commit:
# "OIE replaced with WOAH all Excels"

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Sheet 1")
$wsRefs = $wb.Worksheets.Item("References")

# Cells on "Sheet 1" (column E) that mention OIE and must become WOAH.
$dataCells = @("E5", "E6", "E7", "E14", "E17", "E31", "E53", "E65", "E77", "E79", "E81", "E93", "E133")

foreach ($addr in $dataCells) {
    $cell = $wsData.Range($addr)
    $text = $cell.Value2
    $cell.Value2 = $text.Replace("OIE", "WOAH")
}

# Cells on "References" (column C) that mention OIE and must become WOAH.
$refCells = @("C2", "C5", "C8", "C9", "C10")

foreach ($addr in $refCells) {
    $cell = $wsRefs.Range($addr)
    $text = $cell.Value2
    $cell.Value2 = $text.Replace("OIE", "WOAH")
}
